# fix some report's table err
# Applies the source-data corrections (column D "Plan" figures) for the
# P&L reconciliation block (rows 69-81), restores the SUM roll-up formulas
# on L69/L73, adds the E69-E73+E80 roll-up on E81, switches the O/Q ratio
# columns to a div/0-safe IF() formula, bumps several row heights to 15,
# and updates the sheet/window view state (zoom + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Corrected "Plan" (column D) source values
# ---------------------------------------------------------------------
$ws.Range("D69").Value = 35127792.6561062
$ws.Range("D70").Value = 10876005
$ws.Range("D71").Value = 23426769.8761062
$ws.Range("D72").Value = 825017.78
$ws.Range("D73").Value = 34621177.9939547
$ws.Range("D74").Value = 33183618.0825148
$ws.Range("D75").Value = 981.006666666667
$ws.Range("D76").Value = 1051656.26224918
$ws.Range("D77").Value = 245218.265555556
$ws.Range("D78").Value = 108177.162801852
$ws.Range("D79").Value = 31527.2141666667
$ws.Range("D80").Value = -15418.4689269971
$ws.Range("D81").Value = 491196.193224533

# ---------------------------------------------------------------------
# Formula fixes / additions
# ---------------------------------------------------------------------
# L69 / L73 were plain numbers but should roll up their children, same
# shape as the mirrored K69 / K73 "Plan" columns.
$ws.Range("L69").Formula = "=SUM(L70:L72)"
$ws.Range("L73").Formula = "=SUM(L74:L79)"

# E81 ("Profit before tax") should be derived, not a hard-coded number.
$ws.Range("E81").Formula = "=E69-E73+E80"

# O69 / Q69 are standalone formulas; O70:O81 / Q70:Q81 are one shared
# formula group each (master cell O70 / Q70) - write the whole range so
# the shared-formula group is preserved instead of being split apart.
$ws.Range("O69").Formula = '=IF(K69<>0,L69/K69,"-")'
$ws.Range("Q69").Formula = '=IF(M69<>0,L69/M69,"-")'
$ws.Range("O70:O81").Formula = '=IF(K70<>0,L70/K70,"-")'
$ws.Range("Q70:Q81").Formula = '=IF(M70<>0,L70/M70,"-")'

# ---------------------------------------------------------------------
# Row height bumps (rows 72 / 75-81 already sit at 15)
# ---------------------------------------------------------------------
$ws.Rows.Item(68).RowHeight = 15
$ws.Rows.Item(69).RowHeight = 15
$ws.Rows.Item(70).RowHeight = 15
$ws.Rows.Item(71).RowHeight = 15
$ws.Rows.Item(73).RowHeight = 15
$ws.Rows.Item(74).RowHeight = 15

# ---------------------------------------------------------------------
# View state: zoom out a bit and move the selection down to the bottom
# of the reconciliation block.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 145
$ws.Range("Q75").Select()
